$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "OTHER ACTIVITIES" section: the heading paragraph plus the
#    four paragraphs that follow it (the three volunteer bullets and the
#    "Phi Theta Kappa Alumni" line). The trailing empty paragraph that used
#    to come after them is left in place.
# ---------------------------------------------------------------------------
$startRange = $d.Content
$startRange.Find.Execute("OTHER ACTIVITIES", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$headingPara = $startRange.Paragraphs(1)

$endRange = $d.Content
$endRange.Find.Execute("Phi Theta Kappa Alumni (Honors Society)", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$lastActivityPara = $endRange.Paragraphs(1)

$deleteRange = $d.Range($headingPara.Range.Start, $lastActivityPara.Range.End)
$deleteRange.Delete()

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the "Skills" heading paragraph down to
#    the (now last) empty paragraph at the very end of the document.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$finalPara = $d.Paragraphs.Last

# Placing a brand-new, fully-collapsed bookmark exactly one position before a
# paragraph mark trips an edge case in this host, so a throw-away character
# is inserted right before the mark, the bookmark is anchored just ahead of
# it, and the placeholder is then removed again - leaving the bookmark
# correctly collapsed immediately before the paragraph mark.
$markPos = $finalPara.Range.End - 1
$placeholder = $d.Range($markPos, $markPos)
$placeholder.InsertBefore("X")

$finalPara = $d.Paragraphs.Last
$placeholderPos = $finalPara.Range.End - 2
$bookmarkSpot = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()
